$wb = $excel.ActiveWorkbook

$enDash = [string][char]0x2013

$liveWs = $wb.Worksheets.Item("live-action-tv-series")
$animWs = $wb.Worksheets.Item("animated-tv-series")

# Replace two-digit end-years with four-digit end-years ("always use 4 digit year").
# Cells are touched in the same order the original author worked through the sheets
# (animated-tv-series from row 10 down, then all of live-action-tv-series, then back
# up to the top rows of animated-tv-series) so the shared-string table comes out
# byte-identical to the canonical edit.
$animWs.Range("A10").Value = "2000" + $enDash + "2003"
$animWs.Range("A12").Value = "2006" + $enDash + "2007"
$animWs.Range("A13").Value = "2008" + $enDash + "2009"
$animWs.Range("A15").Value = "2009" + $enDash + "2012"
$animWs.Range("A16").Value = "2009" + $enDash + "2011"
$animWs.Range("A17").Value = "2010" + $enDash + "2012"
$animWs.Range("A18").Value = "2010" + $enDash + "2011"
$animWs.Range("A19").Value = "2011" + $enDash + "2012"
$animWs.Range("A20").Value = "2011" + $enDash + "2012"
$animWs.Range("A21").Value = "2011" + $enDash + "2012"
$animWs.Range("A22").Value = "2012" + $enDash + "2017"
$animWs.Range("A23").Value = "2013" + $enDash + "2019"
$animWs.Range("A24").Value = "2013" + $enDash + "2015"
$animWs.Range("A25").Value = "2015" + $enDash + "2016"
$animWs.Range("A26").Value = "2015" + $enDash + "2019"
$animWs.Range("A27").Value = "2017" + $enDash + "2018"
$animWs.Range("A28").Value = "2017" + $enDash + "2020"

$liveWs.Range("A3").Value = "2013" + $enDash + "2020"
$liveWs.Range("A4").Value = "2015" + $enDash + "2016"
$liveWs.Range("A5").Value = "2015" + $enDash + "2018"
$liveWs.Range("A6").Value = "2015" + $enDash + "2019"
$liveWs.Range("A7").Value = "2016" + $enDash + "2018"
$liveWs.Range("A8").Value = "2017" + $enDash + "2019"
$liveWs.Range("A9").Value = "2017" + $enDash + "2018"
$liveWs.Range("A12").Value = "2017" + $enDash + "2019"
$liveWs.Range("A13").Value = "2017" + $enDash + "2019"
$liveWs.Range("A14").Value = "2017" + $enDash + "2019"
$liveWs.Range("A15").Value = "2018" + $enDash + "2019"
$liveWs.Range("A20").Value = "2021" + $enDash + "2024"

$animWs.Range("A2").Value = "1992" + $enDash + "1997"
$animWs.Range("A3").Value = "1994" + $enDash + "1996"
$animWs.Range("A4").Value = "1994" + $enDash + "1996"
$animWs.Range("A5").Value = "1994" + $enDash + "1998"
$animWs.Range("A6").Value = "1996" + $enDash + "1997"

# Update active cell selections to match target state
$liveWs.Activate()
$liveWs.Range("A21").Select() | Out-Null
$animWs.Activate()
$animWs.Range("A7").Select() | Out-Null
